# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and bumps the
# related timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and generate date (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 15:04:55"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 15:04:51"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 15:04:55"

# Column widths grow to fit the longer "Ready for handoff" text.
# (16.3 is the closest COM-settable ColumnWidth to the target stored
# OOXML column width of ~17.216.)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
